$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark these videos (in the "Video" column, G) as done - "d" - now that
# they've been processed.
$ws.Range("G36").Value = "d"
$ws.Range("G41").Value = "d"
$ws.Range("G42").Value = "d"
$ws.Range("G43").Value = "d"
$ws.Range("G44").Value = "d"

# Leave the view scrolled/selected where the editor ended up working.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("G45").Select()
